$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking scrape stores numeric-looking prices/percent deltas as literal
# text (e.g. "1.00", "9.10", thousand-dot "76.255.17"). Excel auto-coerces a
# bare numeric-looking string typed into Value to a real number, which would
# lose the text semantics/trailing zeros. Force text by prepending the classic
# apostrophe quote-prefix, then restore the cell's style to Normal so we do not
# leave a stray quote-prefixed / text-formatted style on the cell.
function Set-TextValue($range, $value) {
    $range.Value = "'" + $value
    $range.Style = 'Normal'
}

$ws.Range('D2').Value = '76.535.96'
$ws.Range('E2').Value = '  +2.24%  '
$ws.Range('D3').Value = '2.872.22'
$ws.Range('E3').Value = '  +7.83%  '
Set-TextValue $ws.Range('D4') '1.00'
$ws.Range('E4').Value = '  -0.08%  '
Set-TextValue $ws.Range('D5') '196.49'
$ws.Range('E5').Value = '  +5.48%  '
Set-TextValue $ws.Range('D6') '599.66'
$ws.Range('E6').Value = '  +2.38%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('E8').Value = '  +3.54%  '
Set-TextValue $ws.Range('D9') '0.196'
$ws.Range('E9').Value = '  +1.55%  '
$ws.Range('D10').Value = '2.877.69'
$ws.Range('E10').Value = '  +7.93%  '
Set-TextValue $ws.Range('D11') '0.393'
$ws.Range('E11').Value = '  +10.37%  '
Set-TextValue $ws.Range('D12') '0.161'
$ws.Range('E12').Value = '  -1.77%  '
Set-TextValue $ws.Range('D13') '4.92'
$ws.Range('E13').Value = '  +4.05%  '
$ws.Range('D14').Value = '3.397.73'
$ws.Range('E14').Value = '  +7.65%  '
$ws.Range('D15').Value = '76.350.37'
$ws.Range('E15').Value = '  +2.17%  '
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue $ws.Range('D16') '0.0000190'
$ws.Range('E16').Value = '  +2.38%  '
$ws.Range('B17').Value = 'Avalanche'
$ws.Range('C17').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue $ws.Range('D17') '27.53'
$ws.Range('E17').Value = '  +3.99%  '
$ws.Range('D18').Value = '2.854.37'
$ws.Range('E18').Value = '  +7.44%  '
Set-TextValue $ws.Range('D19') '9.10'
Set-TextValue $ws.Range('D20') '12.50'
$ws.Range('E20').Value = '  +5.53%  '
Set-TextValue $ws.Range('D21') '385.22'
$ws.Range('E21').Value = '  +3.89%  '
$ws.Range('E22').Value = '  +4.36%  '
Set-TextValue $ws.Range('D23') '4.15'
$ws.Range('E23').Value = '  +1.76%  '
Set-TextValue $ws.Range('D24') '71.87'
$ws.Range('E24').Value = '  +3.37%  '
$ws.Range('E25').Value = '  -0.02%  '
$ws.Range('D26').Value = '3.008.42'
Set-TextValue $ws.Range('D27') '4.24'
$ws.Range('E27').Value = '  +2.46%  '
Set-TextValue $ws.Range('D28') '9.80'
$ws.Range('E28').Value = '  +5.01%  '
Set-TextValue $ws.Range('D29') '0.0000106'
$ws.Range('E29').Value = '  +11.83%  '
$ws.Range('E30').Value = '  -0.01%  '
$ws.Range('E31').Value = '  +0.07%  '
Set-TextValue $ws.Range('D32') '515.31'
$ws.Range('E32').Value = '  -1.29%  '
Set-TextValue $ws.Range('D33') '7.75'
$ws.Range('E33').Value = '  +1.84%  '
$ws.Range('E34').Value = '  +4.18%  '
Set-TextValue $ws.Range('D35') '0.999'
$ws.Range('E35').Value = '  -0.11%  '
Set-TextValue $ws.Range('D36') '166.87'
$ws.Range('E36').Value = '  +2.03%  '
Set-TextValue $ws.Range('D37') '20.09'
$ws.Range('E37').Value = '  +4.66%  '
$ws.Range('E38').Value = '  +0.92%  '
Set-TextValue $ws.Range('D39') '19.51'
$ws.Range('E39').Value = '  +0.92%  '
Set-TextValue $ws.Range('D40') '186.77'
$ws.Range('E40').Value = '  +9.78%  '
$ws.Range('E41').Value = '  -0.08%  '
$ws.Range('B42').Value = 'PolygonEcosystemToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
Set-TextValue $ws.Range('D42') '0.346'
$ws.Range('E42').Value = '  +5.38%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
Set-TextValue $ws.Range('D43') '5.10'
$ws.Range('E43').Value = '  +2.23%  '
$ws.Range('E44').Value = '  +0.68%  '
$ws.Range('B45').Value = 'ImmutableX'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Range('D45') '1.24'
$ws.Range('E45').Value = '  +4.54%  '
$ws.Range('B46').Value = 'Cronos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws.Range('D46') '0.0914'
$ws.Range('E46').Value = '  +8.94%  '
Set-TextValue $ws.Range('D47') '40.28'
$ws.Range('E47').Value = '  +3.18%  '
$ws.Range('E48').Value = '  +1.78%  '
Set-TextValue $ws.Range('D49') '0.580'
$ws.Range('E49').Value = '  +9.65%  '
Set-TextValue $ws.Range('D50') '0.674'
$ws.Range('E50').Value = '  +14.02%  '
Set-TextValue $ws.Range('D51') '3.75'
$ws.Range('E51').Value = '  +3.32%  '
